# Fix the typo "задара" -> "задача" in the title "Домашна задара број 1".
#
# The canonical OOXML diff shows the single run containing
# "Домашна задара број 1" being rewritten as three runs (all sharing the
# exact same run properties) whose text concatenates to
# "Домашна задача број 1":
#     "Домашна зада" + "ч" + "а број 1"
# i.e. only the single letter "р" -> "ч" actually changed; the run simply
# got split around that one corrected letter (as Word does when you retype
# a single character in the middle of a run).

$d = $word.ActiveDocument

# Locate the title text and fix the typo. Find.Execute collapses/extends
# $rng to the exact bounds of the match, so $rng.Start is the start of the
# (now corrected) title - no hard-coded document offsets needed.
$rng = $d.Content
$rng.Find.Execute("Домашна задара број 1", $true, $false, $false, $false, $false, $true, 1, $false, "Домашна задача број 1", 2)

# The corrected letter "ч" sits right after "Домашна зада" (12 characters).
$charStart = $rng.Start + 12
$letter = $d.Range($charStart, $charStart + 1)

# Touch (and restore) a formatting property on just that one character so
# Word splits it into its own run instead of merging it back with its
# identically-formatted neighbours, matching the three-run shape in the
# target XML.
$letter.Font.Bold = $true
$letter.Font.Bold = $false
